# Nueva actualización V2.2, arreglos en el guardado de precios
#
# Updates the prices on the "estudiante" worksheet: a new column G is added
# (shifting/introducing new price tiers) and several existing prices are
# corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("estudiante")

# Row 1 (header / "cantidad" row) — style s="1" already applied to B1:F1
$ws.Range("C1").Value = 6
$ws.Range("E1").Value = 30
$ws.Range("F1").Value = 50

# New column G1 must inherit the same header formatting as F1 (bold, border,
# centered) before setting its value.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = 100

# Row 2 ("simple")
$ws.Range("B2").Value = 80
$ws.Range("C2").Value = 70
$ws.Range("D2").Value = 60
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 40

# Row 3 ("doble")
$ws.Range("C3").Value = 100
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 70
